$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  A = "Mark";     B = "Mark Dillon";              C = "Dechert LLP"; D = "Ireland";            G = "mark.dillon@dechert.com" },
    @{ Row = 6;  A = "Jay";      B = "Jay Jurata";                C = "Dechert LLP"; D = "Not Found";          G = "jay.jurata@dechert.com" },
    @{ Row = 12; A = "Olaf";     B = "Olaf Fasshauer";            C = "Dechert LLP"; D = "Germany";            G = "olaf.fasshauer@dechert.com" },
    @{ Row = 16; A = "Arne";     B = "Arne Bolch";                C = "Dechert LLP"; D = "Luxembourg";         G = "arne.bolch@dechert.com" },
    @{ Row = 17; A = "Eric";     B = "G Eric Brunstad Jr";        C = "Dechert LLP"; D = "Not Found";          G = "eric.brunstad@dechert.com" },
    @{ Row = 20; A = "Amanjit";  B = "Amanjit K Fagura";          C = "Dechert LLP"; D = "the UAE";             G = "amanjit.fagura@dechert.com" },
    @{ Row = 23; A = "Dean";     B = "Dean Collins";              C = "Dechert LLP"; D = "Singapore";          G = "dean.collins@dechert.com" },
    @{ Row = 31; A = "Eric";     B = "Eric Deltour";              C = "Dechert LLP"; D = "Belgium";            G = "eric.deltour@dechert.com" },
    @{ Row = 37; A = "Daniel";   B = "Daniel Margulies";          C = "Dechert LLP"; D = "Hong Kong";          G = "daniel.margulies@dechert.com" },
    @{ Row = 40; A = "Karen";    B = "Karen L Anderberg";         C = "Dechert LLP"; D = "England";            G = "karen.anderberg@dechert.com" },
    @{ Row = 45; A = "Stephen";  B = "Stephen D Zide";            C = "Dechert LLP"; D = "Not Found";          G = "stephen.zide@dechert.com" },
    @{ Row = 50; A = "Olivia";   B = "Olivia Bernardeau-Paupe";   C = "Dechert LLP"; D = "France";             G = "olivia.bernardeaupaupe@dechert.com" },
    @{ Row = 51; A = "Amanda";   B = "Amanda K Antons Ph D";      C = "Dechert LLP"; D = "Dominican Republic"; G = "amanda.antons@dechert.com" }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("G$r").Value = $entry.G
}
